# Initial commit of light_source_power
#
# This script transforms the existing "UserExperiment*" sheet family into the
# new "LightSourcePower*" sheet family (reusing their positions/sheetIds),
# then re-creates a fresh "UserExperiment*" sheet family (identical to the
# original content) further down, and appends several brand-new support
# sheets (PowerSample, LightSource, PowerMeter, LightSourcePowerKeyMeasurements,
# SimpleLightSourcePowerKeyMeasurements) plus a trailing "LightSourcePower" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: repurpose the 6 existing UserExperiment* sheets in place.
# ---------------------------------------------------------------------------

# UserExperiment -> LightSourcePowerDataset
#   (new content == old UserExperimentDataset header row)
$ws = $wb.Worksheets.Item("UserExperiment")
$ws.Cells.Clear()
$ws.Name = "LightSourcePowerDataset"
$ws.Range("A1").Value = "input_data"
$ws.Range("B1").Value = "input_parameters"
$ws.Range("C1").Value = "output"
$ws.Range("D1").Value = "sample"
$ws.Range("E1").Value = "microscope"
$ws.Range("F1").Value = "experimenter"
$ws.Range("G1").Value = "acquisition_datetime"
$ws.Range("H1").Value = "processed"
$ws.Range("I1").Value = "data_reference"
$ws.Range("J1").Value = "linked_references"
$ws.Range("K1").Value = "name"
$ws.Range("L1").Value = "description"

# UserExperimentDataset -> SimpleLightSourcePowerDataset
$ws = $wb.Worksheets.Item("UserExperimentDataset")
$ws.Cells.Clear()
$ws.Name = "SimpleLightSourcePowerDataset"
$ws.Range("A1").Value = "measurement_device"
$ws.Range("B1").Value = "power_samples"

# UserExperimentInputData -> LightSourcePowerInputData (now empty)
$ws = $wb.Worksheets.Item("UserExperimentInputData")
$ws.Cells.Clear()
$ws.Name = "LightSourcePowerInputData"

# UserExperimentInputParameters -> LightSourcePowerInputParameters
$ws = $wb.Worksheets.Item("UserExperimentInputParameters")
$ws.Cells.Clear()
$ws.Name = "LightSourcePowerInputParameters"
$ws.Range("A1").Value = "key_measurements"
$ws.Range("B1").Value = "processing_application"
$ws.Range("C1").Value = "processing_version"
$ws.Range("D1").Value = "processing_entity"
$ws.Range("E1").Value = "processing_datetime"
$ws.Range("F1").Value = "processing_log"
$ws.Range("G1").Value = "warnings"
$ws.Range("H1").Value = "errors"
$ws.Range("I1").Value = "validated"
$ws.Range("J1").Value = "validation_datetime"
$ws.Range("K1").Value = "comment"

# UserExperimentOutput -> LightSourcePowerOutput
$ws = $wb.Worksheets.Item("UserExperimentOutput")
$ws.Cells.Clear()
$ws.Name = "LightSourcePowerOutput"
$ws.Range("A1").Value = "key_measurements"
$ws.Range("B1").Value = "processing_application"
$ws.Range("C1").Value = "processing_version"
$ws.Range("D1").Value = "processing_entity"
$ws.Range("E1").Value = "processing_datetime"
$ws.Range("F1").Value = "processing_log"
$ws.Range("G1").Value = "warnings"
$ws.Range("H1").Value = "errors"
$ws.Range("I1").Value = "validated"
$ws.Range("J1").Value = "validation_datetime"
$ws.Range("K1").Value = "comment"

# UserExperimentKeyMeasurements -> SimpleLightSourcePowerOutput
$ws = $wb.Worksheets.Item("UserExperimentKeyMeasurements")
$ws.Cells.Clear()
$ws.Name = "SimpleLightSourcePowerOutput"
$ws.Range("A1").Value = "key_measurements"
$ws.Range("B1").Value = "processing_application"
$ws.Range("C1").Value = "processing_version"
$ws.Range("D1").Value = "processing_entity"
$ws.Range("E1").Value = "processing_datetime"
$ws.Range("F1").Value = "processing_log"
$ws.Range("G1").Value = "warnings"
$ws.Range("H1").Value = "errors"
$ws.Range("I1").Value = "validated"
$ws.Range("J1").Value = "validation_datetime"
$ws.Range("K1").Value = "comment"

# ---------------------------------------------------------------------------
# Step 2: append brand-new sheets after SimpleLightSourcePowerOutput.
# ---------------------------------------------------------------------------

# LightSourcePowerKeyMeasurements
$after = $wb.Worksheets.Item("SimpleLightSourcePowerOutput")
$ws = $wb.Worksheets.Add($null, $after)
$ws.Name = "LightSourcePowerKeyMeasurements"
$ws.Range("A1").Value = "light_source"
$ws.Range("B1").Value = "power_mean_mw"
$ws.Range("C1").Value = "power_median_mw"
$ws.Range("D1").Value = "power_std_mw"
$ws.Range("E1").Value = "power_min_mw"
$ws.Range("F1").Value = "power_max_mw"
$ws.Range("G1").Value = "linearity"
$ws.Range("H1").Value = "table_data"
$ws.Range("I1").Value = "data_reference"
$ws.Range("J1").Value = "linked_references"
$ws.Range("K1").Value = "name"
$ws.Range("L1").Value = "description"

# SimpleLightSourcePowerKeyMeasurements
$after = $wb.Worksheets.Item("LightSourcePowerKeyMeasurements")
$ws = $wb.Worksheets.Add($null, $after)
$ws.Name = "SimpleLightSourcePowerKeyMeasurements"
$ws.Range("A1").Value = "light_source"
$ws.Range("B1").Value = "power_mw"
$ws.Range("C1").Value = "table_data"
$ws.Range("D1").Value = "data_reference"
$ws.Range("E1").Value = "linked_references"
$ws.Range("F1").Value = "name"
$ws.Range("G1").Value = "description"

# PowerSample
$after = $wb.Worksheets.Item("SimpleLightSourcePowerKeyMeasurements")
$ws = $wb.Worksheets.Add($null, $after)
$ws.Name = "PowerSample"
$ws.Range("A1").Value = "light_source"
$ws.Range("B1").Value = "sampling_datetime"
$ws.Range("C1").Value = "power_mw"

# LightSource
$after = $wb.Worksheets.Item("PowerSample")
$ws = $wb.Worksheets.Add($null, $after)
$ws.Name = "LightSource"
$ws.Range("A1").Value = "wavelength_nm"

# PowerMeter
$after = $wb.Worksheets.Item("LightSource")
$ws = $wb.Worksheets.Add($null, $after)
$ws.Name = "PowerMeter"
$ws.Range("A1").Value = "manufacturer"
$ws.Range("B1").Value = "model"
$ws.Range("C1").Value = "name"
$ws.Range("D1").Value = "description"

# ---------------------------------------------------------------------------
# Step 3: re-create the UserExperiment* sheet family (same content as the
# original, pre-edit sheets) right after PowerMeter.
# ---------------------------------------------------------------------------

# UserExperiment
$after = $wb.Worksheets.Item("PowerMeter")
$ws = $wb.Worksheets.Add($null, $after)
$ws.Name = "UserExperiment"
$ws.Range("A1").Value = "protocol"
$ws.Range("B1").Value = "manufacturer"
$ws.Range("C1").Value = "name"
$ws.Range("D1").Value = "description"

# UserExperimentDataset
$after = $wb.Worksheets.Item("UserExperiment")
$ws = $wb.Worksheets.Add($null, $after)
$ws.Name = "UserExperimentDataset"
$ws.Range("A1").Value = "input_data"
$ws.Range("B1").Value = "input_parameters"
$ws.Range("C1").Value = "output"
$ws.Range("D1").Value = "sample"
$ws.Range("E1").Value = "microscope"
$ws.Range("F1").Value = "experimenter"
$ws.Range("G1").Value = "acquisition_datetime"
$ws.Range("H1").Value = "processed"
$ws.Range("I1").Value = "data_reference"
$ws.Range("J1").Value = "linked_references"
$ws.Range("K1").Value = "name"
$ws.Range("L1").Value = "description"

# UserExperimentInputData
$after = $wb.Worksheets.Item("UserExperimentDataset")
$ws = $wb.Worksheets.Add($null, $after)
$ws.Name = "UserExperimentInputData"
$ws.Range("A1").Value = "user_experiment_images"
$ws.Range("B1").Value = "orthogonal_rois"
$ws.Range("C1").Value = "profile_rois"

# UserExperimentInputParameters
$after = $wb.Worksheets.Item("UserExperimentInputData")
$ws = $wb.Worksheets.Add($null, $after)
$ws.Name = "UserExperimentInputParameters"
$ws.Range("A1").Value = "bit_depth"
$ws.Range("B1").Value = "saturation_threshold"

# UserExperimentOutput
$after = $wb.Worksheets.Item("UserExperimentInputParameters")
$ws = $wb.Worksheets.Add($null, $after)
$ws.Name = "UserExperimentOutput"
$ws.Range("A1").Value = "intensity_profiles"
$ws.Range("B1").Value = "orthogonal_images"
$ws.Range("C1").Value = "fft_images"
$ws.Range("D1").Value = "key_measurements"
$ws.Range("E1").Value = "processing_application"
$ws.Range("F1").Value = "processing_version"
$ws.Range("G1").Value = "processing_entity"
$ws.Range("H1").Value = "processing_datetime"
$ws.Range("I1").Value = "processing_log"
$ws.Range("J1").Value = "warnings"
$ws.Range("K1").Value = "errors"
$ws.Range("L1").Value = "validated"
$ws.Range("M1").Value = "validation_datetime"
$ws.Range("N1").Value = "comment"

# UserExperimentKeyMeasurements
$after = $wb.Worksheets.Item("UserExperimentOutput")
$ws = $wb.Worksheets.Add($null, $after)
$ws.Name = "UserExperimentKeyMeasurements"
$ws.Range("A1").Value = "channel_name"
$ws.Range("B1").Value = "channel_nr"
$ws.Range("C1").Value = "variation_coefficient"
$ws.Range("D1").Value = "saturated_channels"
$ws.Range("E1").Value = "table_data"
$ws.Range("F1").Value = "data_reference"
$ws.Range("G1").Value = "linked_references"
$ws.Range("H1").Value = "name"
$ws.Range("I1").Value = "description"

# ---------------------------------------------------------------------------
# Step 4: append the final, brand-new "LightSourcePower" sheet at the very end.
# ---------------------------------------------------------------------------

$after = $wb.Worksheets.Item("UserExperimentKeyMeasurements")
$ws = $wb.Worksheets.Add($null, $after)
$ws.Name = "LightSourcePower"
$ws.Range("A1").Value = "protocol"
$ws.Range("B1").Value = "manufacturer"
$ws.Range("C1").Value = "name"
$ws.Range("D1").Value = "description"
